# Repull data, push all data, mean calculation
# Update column F (dSF) values for specific rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = -12
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = -5
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = -6
$ws.Range("F11").Value = 1
